$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 <may> -> <made>; C2 15 -> 14
$ws.Range("B2").Value = "<made>"
$ws.Range("C2").Value = 14

# Row 3: C3 18 -> 14
$ws.Range("C3").Value = 14

# Row 4: C4 13 -> 12
$ws.Range("C4").Value = 12

# Row 5: B5 <an> -> <ab>; C5 16 -> 18
$ws.Range("B5").Value = "<ab>"
$ws.Range("C5").Value = 18

# Row 8: C8 17 -> 16
$ws.Range("C8").Value = 16

# Row 9: C9 21 -> 22
$ws.Range("C9").Value = 22

# Row 10: C10 15 -> 18
$ws.Range("C10").Value = 18

# Row 12: C12 18 -> 16
$ws.Range("C12").Value = 16

# Row 13: C13 13 -> 14
$ws.Range("C13").Value = 14

# Row 14: C14 16 -> 15
$ws.Range("C14").Value = 15

# Row 15: C15 20 -> 19
$ws.Range("C15").Value = 19

# Row 16: B16 <with> -> <we>; C16 16 -> 15
$ws.Range("B16").Value = "<we>"
$ws.Range("C16").Value = 15

# Row 18: C18 12 -> 13
$ws.Range("C18").Value = 13
